# Apply the two changes described by the commit:
#   1. Slide 6's table switches to a different table style
#      ({A572A1B1-63A1-4861-ADD3-BEE209893BBB} -> {C5B04F22-B0F5-48C8-ABA6-6066817AEFA1}).
#   2. The deck's theme palette is swapped from the custom "Integral" colour
#      scheme to the stock Office colour scheme (theme1.xml gains the Office
#      palette; the values below are exactly the 12 clrScheme colours from
#      the target theme, expressed as VBA/COM BGR-packed RGB() integers).

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 6 -------------------------------------------
$slide6 = $p.Slides.Item(6)
$tableShape = $null
for ($i = 1; $i -le $slide6.Shapes.Count; $i++) {
    $candidate = $slide6.Shapes.Item($i)
    if ($candidate.HasTable) {
        $tableShape = $candidate
        break
    }
}
if ($tableShape -eq $null) {
    $tableShape = $slide6.Shapes.Item("Google Shape;127;p18")
}
$tableShape.Table.ApplyStyle("{C5B04F22-B0F5-48C8-ABA6-6066817AEFA1}")

# --- 2. Theme colour palette ----------------------------------------------
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink (Office theme values),
# each packed as 0xBBGGRR for the RGB() COM convention.
$officeColors = @(
    0x000000,  # dk1      000000
    0xFFFFFF,  # lt1      FFFFFF
    0x6A5444,  # dk2      44546A
    0xE6E6E7,  # lt2      E7E6E6
    0xD59B5B,  # accent1  5B9BD5
    0x317DED,  # accent2  ED7D31
    0xA5A5A5,  # accent3  A5A5A5
    0x00C0FF,  # accent4  FFC000
    0xC47244,  # accent5  4472C4
    0x47AD70,  # accent6  70AD47
    0xC16305,  # hlink    0563C1
    0x724F95   # folHlink 954F72
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = $officeColors[$i - 1]
}
